$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Header row: add date / legislator_name / legislator_id columns ---
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

$hdr = $ws.Range("H1:J1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# --- Data rows 2-15: fill in the new columns ---
$dateCol = $ws.Range("H2:H15")
$dateCol.NumberFormat = "@"
for ($r = 2; $r -le 15; $r++) {
    $ws.Range("H" + $r).Value = "2012-04-26"
    $ws.Range("I" + $r).Value = "王進士"
    $ws.Range("J" + $r).Value = 1701
}
$dateCol.Style = $ws.Range("G2:G15").Style()

Write-Host "done"
